$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 550
$ws.Range("I2").Value = 800
$ws.Range("K2").Value = 800
$ws.Range("M2").Value = -687
$ws.Range("H32").Value = 301.58334
$ws.Range("I32").Value = 308.75
$ws.Range("J32").Value = 287.25
$ws.Range("K32").Value = 308.75
$ws.Range("L32").Value = 287.25
$ws.Range("M32").Value = 17.25
$ws.Range("N32").Value = -939.25
$ws.Range("H107").Value = 483.69565
$ws.Range("I107").Value = 543.6316
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 543.6316
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1376.3684
$ws.Range("N107").Value = -4039
$ws.Range("H132").Value = 2940.0967
$ws.Range("I132").Value = 3028.1072
$ws.Range("J132").Value = 2118.6667
$ws.Range("K132").Value = 9084.321599999999
$ws.Range("L132").Value = 6356.000100000001
$ws.Range("M132").Value = -6554.321599999999
$ws.Range("N132").Value = -11416.0001
$ws.Range("H135").Value = 13517349
$ws.Range("I135").Value = 704.4516
$ws.Range("K135").Value = 6340.0644
$ws.Range("M135").Value = -3805.0644
$ws.Range("H137").Value = 1641.35
$ws.Range("I137").Value = 1275.1818
$ws.Range("K137").Value = 3825.5454
$ws.Range("M137").Value = -1275.5454
$ws.Range("H138").Value = 12502760
$ws.Range("I138").Value = 47620596
$ws.Range("J138").Value = 3191.4915
$ws.Range("K138").Value = 142861788
$ws.Range("L138").Value = 9574.4745
$ws.Range("M138").Value = -142856648
$ws.Range("N138").Value = -19854.4745
$ws.Range("H141").Value = 1448.3715
$ws.Range("I141").Value = 1020.31036
$ws.Range("J141").Value = 3517.3333
$ws.Range("K141").Value = 3060.93108
$ws.Range("L141").Value = 10551.9999
$ws.Range("M141").Value = 2119.06892
$ws.Range("N141").Value = -20911.9999

# Sheet 2: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4341.5
$ws.Range("I32").Value = 4178.6387
$ws.Range("K32").Value = 4178.6387
$ws.Range("M32").Value = -3891.6387
$ws.Range("H45").Value = 2595.5557
$ws.Range("I45").Value = 2391.2
$ws.Range("J45").Value = 3060
$ws.Range("K45").Value = 2391.2
$ws.Range("L45").Value = 3060
$ws.Range("M45").Value = -2014.2
$ws.Range("N45").Value = -3814
$ws.Range("H61").Value = 440088.94
$ws.Range("I61").Value = 546593.75
$ws.Range("K61").Value = 546593.75
$ws.Range("M61").Value = -546381.75
$ws.Range("H63").Value = 2120.6
$ws.Range("I63").Value = 1801.3334
$ws.Range("K63").Value = 1801.3334
$ws.Range("M63").Value = -1115.3334
$ws.Range("H66").Value = 2120.6
$ws.Range("I66").Value = 1801.3334
$ws.Range("K66").Value = 9006.666999999999
$ws.Range("M66").Value = -5574.666999999999
$ws.Range("H74").Value = 28573250
$ws.Range("I74").Value = 31251918
$ws.Range("K74").Value = 31251918
$ws.Range("M74").Value = -31251044
$ws.Range("H77").Value = 28573250
$ws.Range("I77").Value = 31251918
$ws.Range("K77").Value = 156259590
$ws.Range("M77").Value = -156255222
$ws.Range("H110").Value = 950.53845
$ws.Range("I110").Value = 921.8889
$ws.Range("K110").Value = 921.8889
$ws.Range("M110").Value = 1123.1111
$ws.Range("H132").Value = 11330.52
$ws.Range("I132").Value = 1307.2683
$ws.Range("J132").Value = 56992
$ws.Range("K132").Value = 3921.8049
$ws.Range("L132").Value = 170976
$ws.Range("M132").Value = -1391.8049
$ws.Range("N132").Value = -176036
$ws.Range("H136").Value = 440088.94
$ws.Range("I136").Value = 546593.75
$ws.Range("K136").Value = 1639781.25
$ws.Range("M136").Value = -1637231.25

# Sheet 3: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 17000
$ws.Range("J46").Value = 17000
$ws.Range("L46").Value = 17000
$ws.Range("N46").Value = -17596
$ws.Range("H94").Value = 1015.04346
$ws.Range("I94").Value = 663.1111
$ws.Range("K94").Value = 663.1111
$ws.Range("M94").Value = -212.1111
$ws.Range("H107").Value = 1223.5
$ws.Range("I107").Value = 561.8182
$ws.Range("J107").Value = 2032.2222
$ws.Range("K107").Value = 561.8182
$ws.Range("L107").Value = 2032.2222
$ws.Range("M107").Value = 1358.1818
$ws.Range("N107").Value = -5872.2222

# Sheet 4: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2364.15
$ws.Range("I31").Value = 1400.841
$ws.Range("J31").Value = 5013.25
$ws.Range("K31").Value = 1400.841
$ws.Range("L31").Value = 5013.25
$ws.Range("M31").Value = -1105.841
$ws.Range("N31").Value = -5603.25
$ws.Range("H34").Value = 2364.15
$ws.Range("I34").Value = 1400.841
$ws.Range("J34").Value = 5013.25
$ws.Range("K34").Value = 1400.841
$ws.Range("L34").Value = 5013.25
$ws.Range("M34").Value = -1198.841
$ws.Range("N34").Value = -5417.25
$ws.Range("H58").Value = 27640.842
$ws.Range("I58").Value = 1428.875
$ws.Range("J58").Value = 167438
$ws.Range("K58").Value = 1428.875
$ws.Range("L58").Value = 167438
$ws.Range("M58").Value = -1225.875
$ws.Range("N58").Value = -167844
$ws.Range("H132").Value = 2009.45
$ws.Range("I132").Value = 1522.1538
$ws.Range("J132").Value = 21014
$ws.Range("K132").Value = 4566.4614
$ws.Range("L132").Value = 63042
$ws.Range("M132").Value = -2036.4614
$ws.Range("N132").Value = -68102
$ws.Range("H134").Value = 957.9048
$ws.Range("I134").Value = 803.7222
$ws.Range("J134").Value = 1883
$ws.Range("K134").Value = 2411.1666
$ws.Range("L134").Value = 5649
$ws.Range("M134").Value = 123.8334
$ws.Range("N134").Value = -10719
$ws.Range("H136").Value = 27640.842
$ws.Range("I136").Value = 1428.875
$ws.Range("J136").Value = 167438
$ws.Range("K136").Value = 4286.625
$ws.Range("L136").Value = 502314
$ws.Range("M136").Value = -1736.625
$ws.Range("N136").Value = -507414

# Sheet 5: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3000429.2
$ws.Range("I4").Value = 463.33334
$ws.Range("J4").Value = 7500378
$ws.Range("K4").Value = 1390.00002
$ws.Range("L4").Value = 22501134
$ws.Range("M4").Value = -1278.00002
$ws.Range("N4").Value = -22501358
$ws.Range("H108").Value = 307.5
$ws.Range("I108").Value = 307.5
$ws.Range("K108").Value = 922.5
$ws.Range("M108").Value = 1957.5
$ws.Range("H118").Value = 41671950
$ws.Range("J118").Value = 8899.571
$ws.Range("L118").Value = 26698.713
$ws.Range("N118").Value = -29184.713
$ws.Range("H119").Value = 4908.3335
$ws.Range("I119").Value = 3900
$ws.Range("K119").Value = 11700
$ws.Range("M119").Value = -6862
$ws.Range("H131").Value = 695.8
$ws.Range("J131").Value = 712.8261
$ws.Range("L131").Value = 2138.4783
$ws.Range("N131").Value = -12218.4783
$ws.Range("H132").Value = 948.25
$ws.Range("J132").Value = 934.6667
$ws.Range("L132").Value = 8412.0003
$ws.Range("N132").Value = -13472.0003

# Sheet 6: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11015.538
$ws.Range("I70").Value = 10773.462
$ws.Range("J70").Value = 11257.615
$ws.Range("K70").Value = 10773.462
$ws.Range("L70").Value = 11257.615
$ws.Range("M70").Value = -10503.462
$ws.Range("N70").Value = -11797.615
$ws.Range("H73").Value = 11015.538
$ws.Range("I73").Value = 10773.462
$ws.Range("J73").Value = 11257.615
$ws.Range("K73").Value = 10773.462
$ws.Range("L73").Value = 11257.615
$ws.Range("M73").Value = -9837.462
$ws.Range("N73").Value = -13129.615
$ws.Range("H97").Value = 1720.5834
$ws.Range("I97").Value = 1949.7
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 1949.7
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = -1453.7
$ws.Range("N97").Value = -1567
$ws.Range("H122").Value = 53334492
$ws.Range("I122").Value = 18519276
$ws.Range("K122").Value = 55557828
$ws.Range("M122").Value = -55555378
$ws.Range("H132").Value = 24506.084
$ws.Range("I132").Value = 3938.4546
$ws.Range("K132").Value = 11815.3638
$ws.Range("M132").Value = -9285.363799999999

# Sheet 7: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 503366.25
$ws.Range("I132").Value = 525216.9399999999
$ws.Range("K132").Value = 1575650.82
$ws.Range("M132").Value = -1573120.82
$ws.Range("H136").Value = 1019.65717
$ws.Range("I136").Value = 954.5862
$ws.Range("K136").Value = 2863.7586
$ws.Range("M136").Value = -313.7586000000001

# Sheet 8: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 497.80768
$ws.Range("I132").Value = 484.0392
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 1452.1176
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = 1077.8824
$ws.Range("N132").Value = -8660
$ws.Range("H136").Value = 17799618
$ws.Range("J136").Value = 5715
$ws.Range("L136").Value = 17145
$ws.Range("N136").Value = -22245
$ws.Range("H138").Value = 46666.668
$ws.Range("J138").Value = 46666.668
$ws.Range("L138").Value = 46666.668
$ws.Range("N138").Value = -56946.668
